$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7374.375
$ws.Range("I40").Value = 9999
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 9999
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -9824
$ws.Range("N40").Value = -3350

$ws.Range("H64").Value = 7450.2354
$ws.Range("I64").Value = 5084.375
$ws.Range("J64").Value = 9553.223
$ws.Range("K64").Value = 5084.375
$ws.Range("L64").Value = 9553.223
$ws.Range("M64").Value = -4836.375
$ws.Range("N64").Value = -10049.223

$ws.Range("H67").Value = 7450.2354
$ws.Range("I67").Value = 5084.375
$ws.Range("J67").Value = 9553.223
$ws.Range("K67").Value = 5084.375
$ws.Range("L67").Value = 9553.223
$ws.Range("M67").Value = -4226.375
$ws.Range("N67").Value = -11269.223

$ws.Range("H87").Value = 63057.832
$ws.Range("J87").Value = 91449
$ws.Range("L87").Value = 91449
$ws.Range("N87").Value = -93945

$ws.Range("H90").Value = 63057.832
$ws.Range("J90").Value = 91449
$ws.Range("L90").Value = 274347
$ws.Range("N90").Value = -286827

$ws.Range("H113").Value = 3840.8
$ws.Range("I113").Value = 9000
$ws.Range("J113").Value = 2551
$ws.Range("K113").Value = 9000
$ws.Range("L113").Value = 2551
$ws.Range("M113").Value = -5746
$ws.Range("N113").Value = -9059

$ws.Range("H116").Value = 52198.2
$ws.Range("I116").Value = 77035
$ws.Range("K116").Value = 77035
$ws.Range("M116").Value = -73593

$ws.Range("H132").Value = 2592.762
$ws.Range("I132").Value = 2552.5
$ws.Range("K132").Value = 7657.5
$ws.Range("M132").Value = -5127.5

$ws.Range("H138").Value = 2513.1702
$ws.Range("I138").Value = 1849.3125
$ws.Range("J138").Value = 2855.8064
$ws.Range("K138").Value = 5547.9375
$ws.Range("L138").Value = 8567.4192
$ws.Range("M138").Value = -407.9375
$ws.Range("N138").Value = -18847.4192

$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 473.42856
$ws.Range("I4").Value = 348.54544
$ws.Range("J4").Value = 931.3333
$ws.Range("K4").Value = 348.54544
$ws.Range("L4").Value = 931.3333
$ws.Range("M4").Value = -232.54544
$ws.Range("N4").Value = -1163.3333

$ws.Range("H5").Value = 310.55554
$ws.Range("I5").Value = 224.66667
$ws.Range("K5").Value = 224.66667
$ws.Range("M5").Value = -112.66667

$ws.Range("H32").Value = 5025.9697
$ws.Range("I32").Value = 5124.4194
$ws.Range("K32").Value = 5124.4194
$ws.Range("M32").Value = -4837.4194

$ws.Range("H61").Value = 3279.2727
$ws.Range("I61").Value = 2008.0555
$ws.Range("K61").Value = 2008.0555
$ws.Range("M61").Value = -1796.0555

$ws.Range("H122").Value = 4025.5334
$ws.Range("I122").Value = 3487.4443
$ws.Range("J122").Value = 4832.6665
$ws.Range("K122").Value = 10462.3329
$ws.Range("L122").Value = 14497.9995
$ws.Range("M122").Value = -8012.332900000001
$ws.Range("N122").Value = -19397.9995

$ws.Range("H136").Value = 3279.2727
$ws.Range("I136").Value = 2008.0555
$ws.Range("K136").Value = 6024.166499999999
$ws.Range("M136").Value = -3474.166499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 310.55554
$ws.Range("I4").Value = 224.66667
$ws.Range("K4").Value = 224.66667
$ws.Range("M4").Value = -109.66667

$ws.Range("H20").Value = 1745.95
$ws.Range("I20").Value = 1869.6666
$ws.Range("J20").Value = 1560.375
$ws.Range("K20").Value = 1869.6666
$ws.Range("L20").Value = 1560.375
$ws.Range("M20").Value = -1622.6666
$ws.Range("N20").Value = -2054.375

$ws.Range("H134").Value = 3658.4878
$ws.Range("I134").Value = 1978.1923
$ws.Range("J134").Value = 6571
$ws.Range("K134").Value = 5934.5769
$ws.Range("L134").Value = 19713
$ws.Range("M134").Value = -3399.5769
$ws.Range("N134").Value = -24783

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H7").Value = 379
$ws.Range("I7").Value = 133.33333
$ws.Range("J7").Value = 747.5
$ws.Range("K7").Value = 133.33333
$ws.Range("L7").Value = 747.5
$ws.Range("M7").Value = -20.33332999999999
$ws.Range("N7").Value = -973.5

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H58").Value = 3198
$ws.Range("I58").Value = 1264.9166
$ws.Range("K58").Value = 1264.9166
$ws.Range("M58").Value = -1061.9166

$ws.Range("H86").Value = 4499.5
$ws.Range("I86").Value = 4499.5
$ws.Range("K86").Value = 4499.5
$ws.Range("M86").Value = -3376.5

$ws.Range("H89").Value = 4499.5
$ws.Range("I89").Value = 4499.5
$ws.Range("K89").Value = 22497.5
$ws.Range("M89").Value = -16881.5

$ws.Range("H99").Value = 5842
$ws.Range("J99").Value = 6171.4546
$ws.Range("L99").Value = 6171.4546
$ws.Range("N99").Value = -9167.454600000001

$ws.Range("H126").Value = 5842
$ws.Range("J126").Value = 6171.4546
$ws.Range("L126").Value = 18514.3638
$ws.Range("N126").Value = -23454.3638

$ws.Range("H136").Value = 3198
$ws.Range("I136").Value = 1264.9166
$ws.Range("K136").Value = 3794.7498
$ws.Range("M136").Value = -1244.7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 168310990
$ws.Range("I4").Value = 202935650
$ws.Range("J4").Value = 12500056
$ws.Range("K4").Value = 608806950
$ws.Range("L4").Value = 37500168
$ws.Range("M4").Value = -608806838
$ws.Range("N4").Value = -37500392

$ws.Range("H33").Value = 148
$ws.Range("I33").Value = 161.5
$ws.Range("J33").Value = 107.5
$ws.Range("K33").Value = 969
$ws.Range("L33").Value = 645
$ws.Range("M33").Value = -686
$ws.Range("N33").Value = -1211

$ws.Range("H131").Value = 1671.2667
$ws.Range("I131").Value = 834.6667
$ws.Range("K131").Value = 2504.0001
$ws.Range("M131").Value = 2535.9999

$ws.Range("H140").Value = 2829.4255
$ws.Range("I140").Value = 1557.3529
$ws.Range("K140").Value = 4672.0587
$ws.Range("M140").Value = 507.9412999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H46").Value = 21195.375
$ws.Range("I46").Value = 15125
$ws.Range("J46").Value = 27265.75
$ws.Range("K46").Value = 15125
$ws.Range("L46").Value = 27265.75
$ws.Range("M46").Value = -14969
$ws.Range("N46").Value = -27577.75

$ws.Range("H122").Value = 5227.645
$ws.Range("I122").Value = 2715.0527
$ws.Range("K122").Value = 8145.158100000001
$ws.Range("M122").Value = -5695.158100000001

$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -79900

$ws.Range("H126").Value = 2887.8845
$ws.Range("I126").Value = 2406
$ws.Range("K126").Value = 7218
$ws.Range("M126").Value = -4748

$ws.Range("H132").Value = 1699.1613
$ws.Range("I132").Value = 1678.3572
$ws.Range("J132").Value = 1893.3334
$ws.Range("K132").Value = 5035.071599999999
$ws.Range("L132").Value = 5680.0002
$ws.Range("M132").Value = -2505.071599999999
$ws.Range("N132").Value = -10740.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1601369.4
$ws.Range("I10").Value = 871
$ws.Range("K10").Value = 871
$ws.Range("M10").Value = -731

$ws.Range("H16").Value = 496.84616
$ws.Range("I16").Value = 529.3
$ws.Range("J16").Value = 388.66666
$ws.Range("K16").Value = 529.3
$ws.Range("L16").Value = 388.66666
$ws.Range("M16").Value = -359.3
$ws.Range("N16").Value = -728.66666

$ws.Range("H46").Value = 1530.5
$ws.Range("I46").Value = 1143.3334
$ws.Range("K46").Value = 1143.3334
$ws.Range("M46").Value = -955.3334

$ws.Range("H93").Value = 3534.3076
$ws.Range("I93").Value = 3937.375
$ws.Range("K93").Value = 3937.375
$ws.Range("M93").Value = -2689.375

$ws.Range("H132").Value = 3312.9714
$ws.Range("I132").Value = 2995.0356
$ws.Range("K132").Value = 8985.106800000001
$ws.Range("M132").Value = -6455.106800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 167090670
$ws.Range("J4").Value = 500012500
$ws.Range("L4").Value = 500012500
$ws.Range("N4").Value = -500012726

$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H132").Value = 4295.853
$ws.Range("I132").Value = 2683.25
$ws.Range("J132").Value = 8166.1
$ws.Range("K132").Value = 8049.75
$ws.Range("L132").Value = 24498.3
$ws.Range("M132").Value = -5519.75
$ws.Range("N132").Value = -29558.3
